$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for H2:H13
$values = @(1,0,0,0,1,1,1,0,0,0,0,0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
